# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp shown in A1
$ws.Range("A1").Value = "Datos actualizados a 24 de Agosto de 2020 a las 09:37"

# Swap country labels: Curazao now listed before Puerto Rico
$ws.Range("A198").Value = "Curazao"
$ws.Range("A199").Value = "Puerto Rico"

# Swap country labels: Timor Oriental now listed before Santa Lucia
$ws.Range("A202").Value = "Timor Oriental"
$ws.Range("A203").Value = "Santa Lucia"

# Update statistics for Rusia (row 7)
$ws.Range("B7").Value = 961493
$ws.Range("C7").Value = 4744
$ws.Range("D7").Value = 773095
$ws.Range("E7").Value = 171950
$ws.Range("G7").Value = 65
$ws.Range("H7").Value = 16448

# Update statistics for Ucrania (row 31)
$ws.Range("B31").Value = 106754
$ws.Range("C31").Value = 1799
$ws.Range("D31").Value = 54524
$ws.Range("E31").Value = 49937
$ws.Range("G31").Value = 22
$ws.Range("H31").Value = 2293

# Update statistics for Singapur (row 49)
$ws.Range("B49").Value = 56404
$ws.Range("C49").Value = 51
$ws.Range("E49").Value = 2213

# Update statistics for Armenia (row 57)
$ws.Range("B57").Value = 42825
$ws.Range("C57").Value = 33
$ws.Range("D57").Value = 36049
$ws.Range("E57").Value = 5922
$ws.Range("G57").Value = 2
$ws.Range("H57").Value = 854

# Update statistics for Australia (row 72)
$ws.Range("D72").Value = 19603
$ws.Range("E72").Value = 4796

# Update statistics for El Salvador (row 73)
$ws.Range("B73").Value = 24811
$ws.Range("C73").Value = 189
$ws.Range("D73").Value = 12492
$ws.Range("E73").Value = 11650

# Update statistics for Hungria (row 109)
$ws.Range("B109").Value = 5191
$ws.Range("C109").Value = 36
$ws.Range("E109").Value = 883

# Update statistics for Curazao (now row 198)
$ws.Range("B198").Value = 43
$ws.Range("C198").Value = 4
$ws.Range("D198").Value = 34
$ws.Range("E198").Value = 8
$ws.Range("H198").Value = 1

# Update statistics for Puerto Rico (now row 199)
$ws.Range("D199").Value = 1
$ws.Range("E199").Value = 36
$ws.Range("H199").Value = 2
